$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cryptocurrency price/volume table update (GitHub Actions data refresh).
# Column D ("Price") values are forced to Text format before assignment
# so that Excel does not reinterpret numeric-looking strings (e.g. "1.027")
# as floating point numbers - matching the original inline-string cells.

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.812.55'
$ws.Range("E2").Value = '  +2.95%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.49'
$ws.Range("E3").Value = '  +2.96%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.027'
$ws.Range("E4").Value = '  +2.26%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.92'
$ws.Range("E5").Value = '  +3.05%  '

# Row 6: USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.034'
$ws.Range("E6").Value = '  +3.06%  '

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5218'
$ws.Range("E7").Value = '  +1.42%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3980'
$ws.Range("E8").Value = '  +3.77%  '

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08422'
$ws.Range("E9").Value = '  +2.36%  '

# Row 10: Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("E10").Value = '  +3.02%  '

# Row 11: WrappedEther -> OKB
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.72'
$ws.Range("E11").Value = '  +2.96%  '

# Row 12: Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.324'
$ws.Range("E12").Value = '  +2.28%  '

# Row 13: Solana -> WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.921.85'
$ws.Range("E13").Value = '  +2.97%  '

# Row 14: Chainlink -> Solana
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.73'
$ws.Range("E14").Value = '  +1.02%  '

# Row 15: BinanceUSD -> Chainlink
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.365'
$ws.Range("E15").Value = '  +1.62%  '

# Row 16: ShibaInu -> BinanceUSD
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.040'
$ws.Range("E16").Value = '  +3.54%  '

# Row 17: Litecoin -> ShibaInu
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001117'
$ws.Range("E17").Value = '  +1.93%  '

# Row 18: TRON -> Litecoin
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.69'
$ws.Range("E18").Value = '  +1.33%  '

# Row 19: Avalanche -> TRON
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06821'
$ws.Range("E19").Value = '  +2.59%  '

# Row 20: Dai -> Avalanche
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.08'
$ws.Range("E20").Value = '  +2.28%  '

# Row 21: Uniswap -> Dai
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.035'
$ws.Range("E21").Value = '  +3.23%  '

# Row 22: WrappedBTC -> Uniswap
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.135'
$ws.Range("E22").Value = '  +2.26%  '

# Row 23: Cosmos -> WrappedBTC
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.857.89'
$ws.Range("E23").Value = '  +3.01%  '

# Row 24: Toncoin -> Cosmos
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("E24").Value = '  +2.38%  '

# Row 25: LEO -> Toncoin
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.294'
$ws.Range("E25").Value = '  +2.23%  '

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.136.48'
$ws.Range("E26").Value = '  +2.96%  '

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.76'
$ws.Range("E27").Value = '  +3.63%  '

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.09'
$ws.Range("E28").Value = '  +3.24%  '

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.487'
$ws.Range("E29").Value = '  -0.83%  '

# Row 30: BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.84'
$ws.Range("E30").Value = '  +2.78%  '

# Row 31: Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1074'
$ws.Range("E31").Value = '  +0.81%  '

# Row 32: ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.057'
$ws.Range("E32").Value = '  +2.79%  '

# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.973'
$ws.Range("E33").Value = '  -0.11%  '

# Row 34: HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.696'
$ws.Range("E34").Value = '  +2.72%  '

# Row 35: FraxShare
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.517'
$ws.Range("E35").Value = '  +1.72%  '

# Row 36: VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02482'
$ws.Range("E36").Value = '  +2.78%  '

# Row 37: Hedera
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06658'
$ws.Range("E37").Value = '  +2.69%  '

# Row 38: Algorand
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2239'
$ws.Range("E38").Value = '  +3.30%  '

# Row 39: TheSandbox
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6593'
$ws.Range("E39").Value = '  +1.19%  '

# Row 40: TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.266'
$ws.Range("E40").Value = '  +3.98%  '

# Row 41: ARBITRUM
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.204'
$ws.Range("E41").Value = '  +0.96%  '

# Row 42: InternetComputer(DFINITY)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.044'
$ws.Range("E42").Value = '  +0.04%  '

# Row 43: Aptos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.17'
$ws.Range("E43").Value = '  +0.38%  '

# Row 44: Decentraland
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6204'
$ws.Range("E44").Value = '  +1.15%  '

# Row 45: EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.28'
$ws.Range("E45").Value = '  +2.36%  '

# Row 46: PancakeSwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.779'
$ws.Range("E46").Value = '  +3.11%  '

# Row 47: WEMIXTOKEN
$ws.Range("E47").Value = '  +2.15%  '

# Row 48: NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.028'
$ws.Range("E48").Value = '  +1.11%  '

# Row 49: EOS
$ws.Range("E49").Value = '  +2.66%  '

# Row 50: Quant
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '123.07'
$ws.Range("E50").Value = '  +2.31%  '

# Row 51: Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06982'
$ws.Range("E51").Value = '  +2.14%  '
